$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("jscosc")

$ws.Range("J1").Value = 34.9360568523407
$ws.Range("J2").Value = 43.98625493049622
$ws.Range("J3").Value = 40.27930974960327
$ws.Range("B4").Value = 2589
$ws.Range("E4").Value = 24
$ws.Range("H4").Value = 99.07264296754251
$ws.Range("I4").Value = 0.01205287713841369
$ws.Range("J4").Value = 36.13392066955566
$ws.Range("B5").Value = 2021
$ws.Range("D5").Value = 2004
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = 22
$ws.Range("G5").Value = 98.91411648568608
$ws.Range("H5").Value = 99.20792079207921
$ws.Range("I5").Value = 0.01874691662555501
$ws.Range("J5").Value = 43.43534731864929
$ws.Range("B6").Value = 1758
$ws.Range("D6").Value = 1754
$ws.Range("F6").Value = 8
$ws.Range("G6").Value = 99.54597048808172
$ws.Range("H6").Value = 99.82925441092772
$ws.Range("I6").Value = 0.006239364719228588
$ws.Range("J6").Value = 33.31949186325073
$ws.Range("J7").Value = 37.95962834358215
$ws.Range("J8").Value = 33.25940561294556
$ws.Range("J9").Value = 36.66025424003601
$ws.Range("B10").Value = 1844
$ws.Range("D10").Value = 1792
$ws.Range("E10").Value = 51
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 99.88851727982163
$ws.Range("H10").Value = 97.23277265328269
$ws.Range("I10").Value = 0.02952646239554317
$ws.Range("J10").Value = 37.58621525764465
$ws.Range("B11").Value = 1862
$ws.Range("D11").Value = 1860
$ws.Range("F11").Value = 18
$ws.Range("G11").Value = 99.04153354632588
$ws.Range("H11").Value = 99.9462654486835
$ws.Range("I11").Value = 0.01011176157530601
$ws.Range("J11").Value = 29.06943345069885
$ws.Range("J12").Value = 32.40294361114502
$ws.Range("J13").Value = 31.56879734992981
$ws.Range("J14").Value = 30.94154095649719
$ws.Range("B15").Value = 2280
$ws.Range("E15").Value = 2
$ws.Range("H15").Value = 99.91224221149628
$ws.Range("I15").Value = 0.000877963125548727
$ws.Range("J15").Value = 32.7888879776001
$ws.Range("B16").Value = 1992
$ws.Range("E16").Value = 5
$ws.Range("H16").Value = 99.74886991461577
$ws.Range("I16").Value = 0.002516356316054353
$ws.Range("J16").Value = 36.21353840827942
$ws.Range("J17").Value = 36.66513395309448
$ws.Range("J18").Value = 35.84201073646545
$ws.Range("J19").Value = 35.50488519668579
$ws.Range("J20").Value = 29.18517446517944
$ws.Range("B21").Value = 2599
$ws.Range("D21").Value = 2597
$ws.Range("F21").Value = 3
$ws.Range("G21").Value = 99.88461538461539
$ws.Range("H21").Value = 99.96150885296382
$ws.Range("I21").Value = 0.001537870049980777
$ws.Range("J21").Value = 34.6830677986145
$ws.Range("B22").Value = 1939
$ws.Range("D22").Value = 1938
$ws.Range("F22").Value = 24
$ws.Range("G22").Value = 98.77675840978593
$ws.Range("I22").Value = 0.01222618441161488
$ws.Range("J22").Value = 40.92652559280396
$ws.Range("J23").Value = 32.79071617126465
$ws.Range("B24").Value = 2936
$ws.Range("D24").Value = 2913
$ws.Range("E24").Value = 22
$ws.Range("F24").Value = 66
$ws.Range("G24").Value = 97.78449144008056
$ws.Range("H24").Value = 99.2504258943782
$ws.Range("I24").Value = 0.02953020134228188
$ws.Range("J24").Value = 40.73287200927734
$ws.Range("B25").Value = 2647
$ws.Range("D25").Value = 2646
$ws.Range("F25").Value = 9
$ws.Range("G25").Value = 99.66101694915254
$ws.Range("I25").Value = 0.00338855421686747
$ws.Range("J25").Value = 34.62625765800476
$ws.Range("B26").Value = 1849
$ws.Range("D26").Value = 1845
$ws.Range("E26").Value = 3
$ws.Range("F26").Value = 14
$ws.Range("G26").Value = 99.24690693921463
$ws.Range("H26").Value = 99.83766233766234
$ws.Range("J26").Value = 36.58816623687744
$ws.Range("B27").Value = 2945
$ws.Range("D27").Value = 2938
$ws.Range("E27").Value = 6
$ws.Range("F27").Value = 16
$ws.Range("G27").Value = 99.4583615436696
$ws.Range("H27").Value = 99.79619565217391
$ws.Range("I27").Value = 0.007445008460236886
$ws.Range("J27").Value = 38.25734901428223
$ws.Range("B28").Value = 3005
$ws.Range("E28").Value = 0
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 35.13791680335999
$ws.Range("B29").Value = 2603
$ws.Range("D29").Value = 2601
$ws.Range("E29").Value = 1
$ws.Range("F29").Value = 48
$ws.Range("G29").Value = 98.18799546998868
$ws.Range("H29").Value = 99.96156802459646
$ws.Range("I29").Value = 0.01849056603773585
$ws.Range("J29").Value = 34.54771399497986
$ws.Range("J30").Value = 33.19395327568054
$ws.Range("J31").Value = 35.1749427318573
$ws.Range("B32").Value = 2261
$ws.Range("E32").Value = 3
$ws.Range("H32").Value = 99.86725663716814
$ws.Range("I32").Value = 0.003094606542882405
$ws.Range("J32").Value = 36.94355177879333
$ws.Range("J33").Value = 40.95211029052734
$ws.Range("J34").Value = 40.08252382278442
$ws.Range("J35").Value = 47.30554366111755
$ws.Range("B36").Value = 2424
$ws.Range("D36").Value = 2414
$ws.Range("F36").Value = 12
$ws.Range("G36").Value = 99.50535861500413
$ws.Range("H36").Value = 99.62855963681386
$ws.Range("I36").Value = 0.00865265760197775
$ws.Range("J36").Value = 37.2197790145874
$ws.Range("B37").Value = 2343
$ws.Range("D37").Value = 2342
$ws.Range("F37").Value = 140
$ws.Range("G37").Value = 94.3593875906527
$ws.Range("I37").Value = 0.05638340716874748
$ws.Range("J37").Value = 40.26085591316223
$ws.Range("J38").Value = 31.6109631061554
$ws.Range("B39").Value = 2054
$ws.Range("D39").Value = 2048
$ws.Range("F39").Value = 4
$ws.Range("G39").Value = 99.80506822612085
$ws.Range("H39").Value = 99.75645396980029
$ws.Range("I39").Value = 0.00438382854359474
$ws.Range("J39").Value = 37.00061416625977
$ws.Range("J40").Value = 37.60756015777588
$ws.Range("J41").Value = 35.63766407966614
$ws.Range("J42").Value = 34.92338347434998
$ws.Range("B43").Value = 3077
$ws.Range("D43").Value = 3070
$ws.Range("F43").Value = 8
$ws.Range("G43").Value = 99.74009096816114
$ws.Range("H43").Value = 99.80494148244473
$ws.Range("I43").Value = 0.004546930821695355
$ws.Range("J43").Value = 42.74222373962402
$ws.Range("J44").Value = 37.61645340919495

$ws.Name = "sddkj"

